# Generate Report for Handoff
#
# Updates the "Latest Handoff Datetime" timestamps on the Overview / de-de /
# zh-cn sheets and sets the handback "Priority" column to "ht" for the
# file-path-mismatch rows on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 12, 13, 14)

# "Latest HO Xliff Generate Date" on the Overview sheet (column G).
foreach ($r in $rows) {
    $ws_overview.Cells.Item($r, 7).Value = "2016-08-28 00:20:08"
}

# "Latest Handoff Datetime" on the de-de sheet (column H) mirrors the
# Overview sheet's generate date.
foreach ($r in $rows) {
    $ws_dede.Cells.Item($r, 8).Value = "2016-08-28 00:20:08"
}

# "Latest Handoff Datetime" on the zh-cn sheet (column H).
foreach ($r in $rows) {
    $ws_zhcn.Cells.Item($r, 8).Value = "2016-08-28 00:19:57"
}

# "Priority" column (E) on zh-cn and de-de now flags the handback-priority
# mismatch with "ht" (handoff type).
foreach ($r in $rows) {
    $ws_zhcn.Cells.Item($r, 5).Value = "ht"
    $ws_dede.Cells.Item($r, 5).Value = "ht"
}
